$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- VENDA FINANCIADA (B5:C14) ----
$ws.Range("C6").Value = 90000
$ws.Range("C7").Value = 10000
$ws.Range("C8").Formula = "=C6-C7"
$ws.Range("C9").Value = 36
$ws.Range("C10").Value = 0.01
$ws.Range("C11").Formula = "=PMT(C10,C9,-C8,0,0)"
$ws.Range("C12").Formula = "=C11*C9"
$ws.Range("C13").Formula = "=C12-C8"

# ---- DESCONTO DE DUPLICATA (E5:F12) ----
$ws.Range("F6").Value = 3000
$ws.Range("F7").Value = 43819
$ws.Range("F8").Value = 43728
$ws.Range("F9").Formula = "=(F7-F8)/30"
$ws.Range("F10").Value = 0.02
$ws.Range("F12").Formula = "=PV(F10,F9,0,-F6,0)"
$ws.Range("F11").Formula = "=F6-F12"

# ---- SIMULADOR DE APLICAÇÃO (H5:I11) ----
$ws.Range("I6").Value = 10000
$ws.Range("I7").Value = 43728
$ws.Range("I8").Value = 44094
$ws.Range("I9").Value = 12
$ws.Range("I10").Value = 0.01
$ws.Range("I11").Formula = "=FV(I10,I9,0,-I6,0)"

# ---- Update the sheet view: scroll so column C is leftmost, select H14 ----
$ws.Activate()
$ws.Range("H14").Select()
$excel.ActiveWindow.ScrollColumn = 3
